$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows at row 306, shifting the existing rows 306:368
# down to 309:371 (and dimension grows from A1:T368 to A1:T371).
$ws.Range("A306:A308").EntireRow.Insert()

# Populate the 3 newly-inserted rows with a new "Tuna" price block
# (same Mercado/Region/Producto metadata as the block that used to sit
# here, but a new Fecha and new Volumen/Precio/Precio $/Kg figures).

# Row 306: Especial
$ws.Cells.Item(306, 1).Value  = 9
$ws.Cells.Item(306, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(306, 3).Value  = "Metropolitana"
$ws.Cells.Item(306, 4).Value  = 44798
$ws.Cells.Item(306, 5).Value  = 13
$ws.Cells.Item(306, 6).Value  = "Fruta"
$ws.Cells.Item(306, 7).Value  = 100107
$ws.Cells.Item(306, 8).Value  = "Otros"
$ws.Cells.Item(306, 9).Value  = 100107011
$ws.Cells.Item(306, 10).Value = "Tuna"
$ws.Cells.Item(306, 11).Value = "Sin especificar"
$ws.Cells.Item(306, 12).Value = "Especial"
$ws.Cells.Item(306, 13).Value = 380
$ws.Cells.Item(306, 14).Value = 25000
$ws.Cells.Item(306, 15).Value = 27000
$ws.Cells.Item(306, 16).Value = 25947
$ws.Cells.Item(306, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(306, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(306, 19).Value = 1442
$ws.Cells.Item(306, 20).Value = 18

# Row 307: Primera
$ws.Cells.Item(307, 1).Value  = 9
$ws.Cells.Item(307, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(307, 3).Value  = "Metropolitana"
$ws.Cells.Item(307, 4).Value  = 44798
$ws.Cells.Item(307, 5).Value  = 13
$ws.Cells.Item(307, 6).Value  = "Fruta"
$ws.Cells.Item(307, 7).Value  = 100107
$ws.Cells.Item(307, 8).Value  = "Otros"
$ws.Cells.Item(307, 9).Value  = 100107011
$ws.Cells.Item(307, 10).Value = "Tuna"
$ws.Cells.Item(307, 11).Value = "Sin especificar"
$ws.Cells.Item(307, 12).Value = "Primera"
$ws.Cells.Item(307, 13).Value = 350
$ws.Cells.Item(307, 14).Value = 20000
$ws.Cells.Item(307, 15).Value = 22000
$ws.Cells.Item(307, 16).Value = 20857
$ws.Cells.Item(307, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(307, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(307, 19).Value = 1159
$ws.Cells.Item(307, 20).Value = 18

# Row 308: Segunda
$ws.Cells.Item(308, 1).Value  = 9
$ws.Cells.Item(308, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(308, 3).Value  = "Metropolitana"
$ws.Cells.Item(308, 4).Value  = 44798
$ws.Cells.Item(308, 5).Value  = 13
$ws.Cells.Item(308, 6).Value  = "Fruta"
$ws.Cells.Item(308, 7).Value  = 100107
$ws.Cells.Item(308, 8).Value  = "Otros"
$ws.Cells.Item(308, 9).Value  = 100107011
$ws.Cells.Item(308, 10).Value = "Tuna"
$ws.Cells.Item(308, 11).Value = "Sin especificar"
$ws.Cells.Item(308, 12).Value = "Segunda"
$ws.Cells.Item(308, 13).Value = 270
$ws.Cells.Item(308, 14).Value = 14000
$ws.Cells.Item(308, 15).Value = 16000
$ws.Cells.Item(308, 16).Value = 15111
$ws.Cells.Item(308, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(308, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(308, 19).Value = 840
$ws.Cells.Item(308, 20).Value = 18
